$d = $word.ActiveDocument

# The document starts with a stray "Fdsfd" paragraph (complete with
# spell-check proofing markers) followed by a blank paragraph, then the
# real "Yolanda Brown" paragraph. Remove the bogus paragraph and the
# blank one entirely (text + paragraph marks), leaving "Yolanda Brown"
# as the document's first paragraph.
$p1 = $d.Paragraphs(1)
$p2 = $d.Paragraphs(2)
$junkRange = $d.Range($p1.Range.Start, $p2.Range.End)
$junkRange.Delete()

# Split "Yolanda Brown" into two runs ("Y" and "olanda Brown") by cutting
# the tail of the word to the clipboard and pasting it right back in
# place. This forces Word to create a fresh run for the pasted text
# instead of silently re-merging it with the untouched "Y" run.
$para = $d.Paragraphs(1)
$full = $para.Range
$start = $full.Start
$tailRange = $d.Range($start + 1, $start + 13)
$tailRange.Cut()
$pastePoint = $d.Range($start + 1, $start + 1)
$pastePoint.Paste()
